$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update timing values in rows 5-7 ---
$ws.Range("B5").Value = 0.0006709098815917969
$ws.Range("B6").Value = 0.000560760498046875
$ws.Range("B7").Value = 0.002531051635742188

# --- Convert tuple-style text to list-style text ---
$ws.Range("A8").Value = "[[2, 0], [1, 3], [0, 3], [0, 2], [0, 0], [2, 1], [1, 2], [1, 0], [1, 1], [0, 1]]"
$ws.Range("A48").Value = "[[2, 1], [1, 1], [1, 2], [0, 2], [0, 3], [1, 3], [0, 1], [2, 2], [1, 0], [2, 0]]"
$ws.Range("A98").Value = "[[1, 1], [0, 2], [0, 1], [0, 0], [1, 2], [1, 0], [2, 0], [0, 3], [1, 3], [2, 1]]"

# --- Insert a new row before the old row 134 ("Movement times") ---
# This shifts rows 134-138 down to 135-139.
$ws.Rows.Item(134).Insert()

# --- Populate the newly inserted row 134 with the move_fidelity stat ---
$ws.Range("A134").Value = "move_fidelity"
$ws.Range("B134").Value = 0.9990298597111551

# --- Update the "total time:" value, now living at row 138 ---
$ws.Range("B138").Value = 0.01520490646362305
